# Update "想去人数" (interested-attendee count) figures in column F
# for the "展览" sheet and the "全部类型" sheet, matching the
# upstream data refresh captured in commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (row -> new F value) ---
$sheet1Updates = @{
    2  = 352
    4  = 1284
    6  = 29
    9  = 148
    10 = 3536
    11 = 141
    14 = 49
    16 = 612
    17 = 103
    18 = 769
    20 = 131
    21 = 58
    22 = 65
    23 = 70
    24 = 2728
    25 = 5231
    29 = 3089
    30 = 293
    31 = 2267
    34 = 86
    35 = 131
    36 = 182
    38 = 35
    39 = 464
    40 = 811
    45 = 492
}

# --- Sheet "全部类型" (row -> new F value) ---
$sheet4Updates = @{
    2  = 352
    4  = 1284
    6  = 29
    9  = 148
    10 = 3536
    11 = 141
    15 = 49
    17 = 612
    18 = 103
    19 = 769
    21 = 131
    22 = 58
    23 = 65
    24 = 70
    25 = 2728
    26 = 5232
    30 = 3089
    31 = 293
    32 = 2267
    35 = 86
    36 = 131
    37 = 182
    39 = 35
    40 = 464
    41 = 811
    46 = 492
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
